$wb = $excel.ActiveWorkbook

# --- Remove the stale/obsolete sheets (old province sheets no longer needed) ---
$null = $wb.Worksheets.Item("TINH TAY NINH").Delete()
$null = $wb.Worksheets.Item("QUAN 2").Delete()
$null = $wb.Worksheets.Item("TINH LAM DONG").Delete()

# --- Rename "TINH PHU YEN" -> "HUYEN VINH THANH" ---
$wb.Worksheets.Item("TINH PHU YEN").Name = "HUYEN VINH THANH"

# --- Fill in the corrected company/tax-code data on "TINH BINH DINH" ---
$ws = $wb.Worksheets.Item("TINH BINH DINH")

# Row 1
$ws.Range("A1").Value = 4101598873
$ws.Range("B1").NumberFormat = "mm-dd-yy"
$ws.Range("B1").Value = "08/06/2021"
$ws.Range("C1").Value = "CÔNG TY TNHH VẬN TẢI TUẤN MINH BÌNH ĐỊNH"
$ws.Range("D1").Value = "Lô 15-16 Khu Đô thị An Phú Thịnh, Phường Đống Đa, Thành phố Quy Nhơn, Tỉnh Bình Định"
$ws.Range("E1").Value = 965972999
$ws.Range("F1").Value = "Nguyễn Quy Khoa"

# Row 2
$ws.Range("A2").Value = 4101598880
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = "08/06/2021"
$ws.Range("C2").Value = "CÔNG TY TNHH THƯƠNG MẠI - TỔNG HỢP BÌNH VƯƠNG"
$ws.Range("D2").Value = "Số 295 Nguyễn Thị Minh Khai, Phường Nguyễn Văn Cừ, Thành phố Quy Nhơn, Tỉnh Bình Định"
$ws.Range("E2").Value = 963555405
$ws.Range("F2").Value = "Bùi Quốc Thắng"
